$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header figures -------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 391014
# Cant. Trabajadores (worker count)
$ws.Range("C13").Value = 8

# --- Make room for two more worker rows -----------------------------------
# Existing data rows are 16-21 (6 workers); insert two fresh rows right
# after row 21 so the trailing "Observaciones / firma" block (old rows
# 26-27) is pushed down to rows 28-29, matching the new layout.
$ws.Rows("22:23").Insert()

# Row 22 should look like a normal (non-last) data row: clone formatting
# from row 16.
$ws.Range("B16:J16").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)

# Row 23 becomes the new "last" data row (double border under-line style):
# clone that formatting from row 21, which currently still has it.
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# Row 21 is no longer the last row, so restyle it like a normal row too.
$ws.Range("B16:J16").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Worker detail rows -----------------------------------------------
# All periods move from 2506 to 2507, a couple of rows are reshuffled and
# two brand-new employees are appended before the closing "last" row.

# Row 16: LUIS EDUARDO MERCADO ROBLES (unchanged person, new period)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1043297908"
$ws.Range("D16").Value = "LUIS EDUARDO MERCADO ROBLES"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Row 17: JESUS ALBERTO RIOS PEREIRA (unchanged person, new period)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1044924108"
$ws.Range("D17").Value = "JESUS ALBERTO RIOS PEREIRA"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: ISABEL MARIA CASTRO CANTILLO
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45556298"
$ws.Range("D18").Value = "ISABEL MARIA CASTRO CANTILLO"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: MATEO DE JESUS MENDOZA GOMEZ
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1041974534"
$ws.Range("D19").Value = "MATEO DE JESUS MENDOZA GOMEZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20: JESSE DE JESUS OSORIO CASTELLON (new worker)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143415630"
$ws.Range("D20").Value = "JESSE DE JESUS OSORIO CASTELLON"
$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 24700
$ws.Range("G20").Value = 1425000

# Row 21: MICHELL MATURANA RUIZ
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1042578779"
$ws.Range("D21").Value = "MICHELL MATURANA RUIZ"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# Row 22: LINDA LUZ NOVOA CANTILLO
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1192717234"
$ws.Range("D22").Value = "LINDA LUZ NOVOA CANTILLO"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

# Row 23: MARIA FERNANDA ESPINOSA PADILLA (new worker, last row)
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1005682282"
$ws.Range("D23").Value = "MARIA FERNANDA ESPINOSA PADILLA"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 24674
$ws.Range("G23").Value = 1423500

# --- Column D grew a bit wider to fit the longer new names ---------------
$ws.Columns("D").ColumnWidth = 34.14

Write-Host "Edit complete"
